$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D, shifting existing D:K to E:L
$ws.Columns("D:D").Insert()

# Copy cell formatting (number format, font) from column E into new column D
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set cell values for columns D through L per the refreshed financial data
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43100
$ws.Range("F7").Value = 42735
$ws.Range("G7").Value = 42369
$ws.Range("H7").Value = 42004
$ws.Range("I7").Value = "NA"
$ws.Range("J7").Value = "NA"
$ws.Range("K7").Value = "NA"
$ws.Range("D8").Value = 3705700
$ws.Range("E8").Value = 3590000
$ws.Range("F8").Value = 2723800
$ws.Range("G8").Value = 1217300
$ws.Range("H8").Value = 1204600
$ws.Range("I8").Value = "NA"
$ws.Range("J8").Value = "NA"
$ws.Range("K8").Value = "NA"
$ws.Range("D9").Value = 889800
$ws.Range("E9").Value = 876200
$ws.Range("F9").Value = 677200
$ws.Range("G9").Value = 337700
$ws.Range("H9").Value = 324300
$ws.Range("I9").Value = "NA"
$ws.Range("J9").Value = "NA"
$ws.Range("K9").Value = "NA"
$ws.Range("D10").Value = 2815900
$ws.Range("E10").Value = 2713800
$ws.Range("F10").Value = 2046600
$ws.Range("G10").Value = 879600
$ws.Range("H10").Value = 880300
$ws.Range("I10").Value = "NA"
$ws.Range("J10").Value = "NA"
$ws.Range("K10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("G12").Value = "NA"
$ws.Range("H12").Value = "NA"
$ws.Range("I12").Value = "NA"
$ws.Range("J12").Value = "NA"
$ws.Range("K12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("D14").Value = 613100
$ws.Range("E14").Value = 759400
$ws.Range("F14").Value = 152900
$ws.Range("G14").Value = 19800
$ws.Range("H14").Value = 31900
$ws.Range("I14").Value = "NA"
$ws.Range("J14").Value = "NA"
$ws.Range("K14").Value = "NA"
$ws.Range("D15").Value = 829800
$ws.Range("E15").Value = 793700
$ws.Range("F15").Value = 587300
$ws.Range("G15").Value = 216400
$ws.Range("H15").Value = 216700
$ws.Range("I15").Value = "NA"
$ws.Range("J15").Value = "NA"
$ws.Range("K15").Value = "NA"
$ws.Range("D17").Value = 3761400
$ws.Range("E17").Value = 3804700
$ws.Range("F17").Value = 2407600
$ws.Range("G17").Value = 969200
$ws.Range("H17").Value = 988000
$ws.Range("I17").Value = "NA"
$ws.Range("J17").Value = "NA"
$ws.Range("K17").Value = "NA"
$ws.Range("D18").Value = -55700
$ws.Range("E18").Value = -214700
$ws.Range("F18").Value = 316200
$ws.Range("G18").Value = 248100
$ws.Range("H18").Value = 216600
$ws.Range("I18").Value = "NA"
$ws.Range("J18").Value = "NA"
$ws.Range("K18").Value = "NA"
$ws.Range("D20").Value = -85300
$ws.Range("E20").Value = -54700
$ws.Range("F20").Value = -99900
$ws.Range("G20").Value = 2100
$ws.Range("H20").Value = -52100
$ws.Range("I20").Value = "NA"
$ws.Range("J20").Value = "NA"
$ws.Range("K20").Value = "NA"
$ws.Range("D21").Value = 688800
$ws.Range("E21").Value = 524300
$ws.Range("F21").Value = 803600
$ws.Range("G21").Value = 466600
$ws.Range("H21").Value = 381200
$ws.Range("I21").Value = "NA"
$ws.Range("J21").Value = "NA"
$ws.Range("K21").Value = "NA"
$ws.Range("D22").Value = 443700
$ws.Range("E22").Value = 381800
$ws.Range("F22").Value = 314400
$ws.Range("G22").Value = 157900
$ws.Range("H22").Value = 140400
$ws.Range("I22").Value = "NA"
$ws.Range("J22").Value = "NA"
$ws.Range("K22").Value = "NA"
$ws.Range("D23").Value = -584700
$ws.Range("E23").Value = -651200
$ws.Range("F23").Value = -98100
$ws.Range("G23").Value = 92300
$ws.Range("H23").Value = 24100
$ws.Range("I23").Value = "NA"
$ws.Range("J23").Value = "NA"
$ws.Range("K23").Value = "NA"
$ws.Range("D24").Value = 51100
$ws.Range("E24").Value = 231200
$ws.Range("F24").Value = 305900
$ws.Range("G24").Value = 46500
$ws.Range("H24").Value = 14400
$ws.Range("I24").Value = "NA"
$ws.Range("J24").Value = "NA"
$ws.Range("K24").Value = "NA"
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("D26").Value = -635800
$ws.Range("E26").Value = -882400
$ws.Range("F26").Value = -404000
$ws.Range("G26").Value = 45800
$ws.Range("H26").Value = 9700
$ws.Range("I26").Value = "NA"
$ws.Range("J26").Value = "NA"
$ws.Range("K26").Value = "NA"
$ws.Range("D27").Value = -345200
$ws.Range("E27").Value = -861800
$ws.Range("F27").Value = -432300
$ws.Range("G27").Value = 38000
$ws.Range("H27").Value = 12000
$ws.Range("I27").Value = "NA"
$ws.Range("J27").Value = "NA"
$ws.Range("K27").Value = "NA"
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = 83700
$ws.Range("F29").Value = "NA"
$ws.Range("G29").Value = "NA"
$ws.Range("H29").Value = "NA"
$ws.Range("I29").Value = "NA"
$ws.Range("J29").Value = "NA"
$ws.Range("K29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("D32").Value = 85300
$ws.Range("E32").Value = 54700
$ws.Range("F32").Value = 99900
$ws.Range("G32").Value = -2100
$ws.Range("H32").Value = 52100
$ws.Range("I32").Value = "NA"
$ws.Range("J32").Value = "NA"
$ws.Range("K32").Value = "NA"
$ws.Range("D33").Value = -345200
$ws.Range("E33").Value = -778100
$ws.Range("F33").Value = -432300
$ws.Range("G33").Value = 38000
$ws.Range("H33").Value = 12000
$ws.Range("I33").Value = "NA"
$ws.Range("J33").Value = "NA"
$ws.Range("K33").Value = "NA"
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("D35").Value = -345200
$ws.Range("E35").Value = -778100
$ws.Range("F35").Value = -432300
$ws.Range("G35").Value = 38000
$ws.Range("H35").Value = 12000
$ws.Range("I35").Value = "NA"
$ws.Range("J35").Value = "NA"
$ws.Range("K35").Value = "NA"
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43100
$ws.Range("F38").Value = 42735
$ws.Range("G38").Value = 42369
$ws.Range("H38").Value = 42004
$ws.Range("I38").Value = "NA"
$ws.Range("J38").Value = "NA"
$ws.Range("K38").Value = "NA"
$ws.Range("D41").Value = 631000
$ws.Range("E41").Value = 529900
$ws.Range("F41").Value = 552600
$ws.Range("G41").Value = 274500
$ws.Range("H41").Value = "NA"
$ws.Range("I41").Value = "NA"
$ws.Range("J41").Value = "NA"
$ws.Range("K41").Value = "NA"
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("D43").Value = 607300
$ws.Range("E43").Value = 556500
$ws.Range("F43").Value = 617800
$ws.Range("G43").Value = 91500
$ws.Range("H43").Value = "NA"
$ws.Range("I43").Value = "NA"
$ws.Range("J43").Value = "NA"
$ws.Range("K43").Value = "NA"
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("D45").Value = 406500
$ws.Range("E45").Value = 288400
$ws.Range("F45").Value = 338800
$ws.Range("G45").Value = 73000
$ws.Range("H45").Value = "NA"
$ws.Range("I45").Value = "NA"
$ws.Range("J45").Value = "NA"
$ws.Range("K45").Value = "NA"
$ws.Range("D46").Value = 1644800
$ws.Range("E46").Value = 1374800
$ws.Range("F46").Value = 1509200
$ws.Range("G46").Value = 439000
$ws.Range("H46").Value = "NA"
$ws.Range("I46").Value = "NA"
$ws.Range("J46").Value = "NA"
$ws.Range("K46").Value = "NA"
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("D48").Value = 4236900
$ws.Range("E48").Value = 4169200
$ws.Range("F48").Value = 3860900
$ws.Range("G48").Value = 843500
$ws.Range("H48").Value = "NA"
$ws.Range("I48").Value = "NA"
$ws.Range("J48").Value = "NA"
$ws.Range("K48").Value = "NA"
$ws.Range("D49").Value = 6861500
$ws.Range("E49").Value = 7555200
$ws.Range("F49").Value = 8195200
$ws.Range("G49").Value = 1498900
$ws.Range("H49").Value = "NA"
$ws.Range("I49").Value = "NA"
$ws.Range("J49").Value = "NA"
$ws.Range("K49").Value = "NA"
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("D52").Value = 703400
$ws.Range("E52").Value = 517700
$ws.Range("F52").Value = 578600
$ws.Range("G52").Value = 456700
$ws.Range("H52").Value = "NA"
$ws.Range("I52").Value = "NA"
$ws.Range("J52").Value = "NA"
$ws.Range("K52").Value = "NA"
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("D54").Value = 13446600
$ws.Range("E54").Value = 13616900
$ws.Range("F54").Value = 14143900
$ws.Range("G54").Value = 3238100
$ws.Range("H54").Value = "NA"
$ws.Range("I54").Value = "NA"
$ws.Range("J54").Value = "NA"
$ws.Range("K54").Value = "NA"
$ws.Range("D57").Value = 297400
$ws.Range("E57").Value = 286800
$ws.Range("F57").Value = 219400
$ws.Range("G57").Value = 54200
$ws.Range("H57").Value = "NA"
$ws.Range("I57").Value = "NA"
$ws.Range("J57").Value = "NA"
$ws.Range("K57").Value = "NA"
$ws.Range("D58").Value = 302500
$ws.Range("E58").Value = 263300
$ws.Range("F58").Value = 150800
$ws.Range("G58").Value = 800
$ws.Range("H58").Value = "NA"
$ws.Range("I58").Value = "NA"
$ws.Range("J58").Value = "NA"
$ws.Range("K58").Value = "NA"
$ws.Range("D59").Value = 1008800
$ws.Range("E59").Value = 1036800
$ws.Range("F59").Value = 977800
$ws.Range("G59").Value = 343500
$ws.Range("H59").Value = "NA"
$ws.Range("I59").Value = "NA"
$ws.Range("J59").Value = "NA"
$ws.Range("K59").Value = "NA"
$ws.Range("D60").Value = 1608700
$ws.Range("E60").Value = 1586900
$ws.Range("F60").Value = 1348000
$ws.Range("G60").Value = 398500
$ws.Range("H60").Value = "NA"
$ws.Range("I60").Value = "NA"
$ws.Range("J60").Value = "NA"
$ws.Range("K60").Value = "NA"
$ws.Range("D61").Value = 6379600
$ws.Range("E61").Value = 6108200
$ws.Range("F61").Value = 5897100
$ws.Range("G61").Value = 2304600
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("D62").Value = 1334900
$ws.Range("E62").Value = 1231200
$ws.Range("F62").Value = 1238400
$ws.Range("G62").Value = 264200
$ws.Range("H62").Value = "NA"
$ws.Range("I62").Value = "NA"
$ws.Range("J62").Value = "NA"
$ws.Range("K62").Value = "NA"
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("D66").Value = 10334000
$ws.Range("E66").Value = 10287300
$ws.Range("F66").Value = 9964300
$ws.Range("G66").Value = 3030600
$ws.Range("H66").Value = "NA"
$ws.Range("I66").Value = "NA"
$ws.Range("J66").Value = "NA"
$ws.Range("K66").Value = "NA"
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("D72").Value = -1367000
$ws.Range("E72").Value = -1010700
$ws.Range("F72").Value = -232600
$ws.Range("G72").Value = 199700
$ws.Range("H72").Value = "NA"
$ws.Range("I72").Value = "NA"
$ws.Range("J72").Value = "NA"
$ws.Range("K72").Value = "NA"
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("D76").Value = 3112600
$ws.Range("E76").Value = 3329600
$ws.Range("F76").Value = 4179600
$ws.Range("G76").Value = 207500
$ws.Range("H76").Value = "NA"
$ws.Range("I76").Value = "NA"
$ws.Range("J76").Value = "NA"
$ws.Range("K76").Value = "NA"
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43100
$ws.Range("F80").Value = 42735
$ws.Range("G80").Value = 42369
$ws.Range("H80").Value = 42004
$ws.Range("I80").Value = "NA"
$ws.Range("J80").Value = "NA"
$ws.Range("K80").Value = "NA"
$ws.Range("D81").Value = -345200
$ws.Range("E81").Value = -778100
$ws.Range("F81").Value = -432300
$ws.Range("G81").Value = 38000
$ws.Range("H81").Value = 12000
$ws.Range("I81").Value = "NA"
$ws.Range("J81").Value = "NA"
$ws.Range("K81").Value = "NA"
$ws.Range("D83").Value = 829800
$ws.Range("E83").Value = 793700
$ws.Range("F83").Value = 587300
$ws.Range("G83").Value = 216400
$ws.Range("H83").Value = 216700
$ws.Range("I83").Value = "NA"
$ws.Range("J83").Value = "NA"
$ws.Range("K83").Value = "NA"
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("D89").Value = 816800
$ws.Range("E89").Value = 573200
$ws.Range("F89").Value = 468200
$ws.Range("G89").Value = 310200
$ws.Range("H89").Value = 289100
$ws.Range("I89").Value = "NA"
$ws.Range("J89").Value = "NA"
$ws.Range("K89").Value = "NA"
$ws.Range("D91").Value = -776400
$ws.Range("E91").Value = -639300
$ws.Range("F91").Value = -490400
$ws.Range("G91").Value = -227200
$ws.Range("H91").Value = -223100
$ws.Range("I91").Value = "NA"
$ws.Range("J91").Value = "NA"
$ws.Range("K91").Value = "NA"
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("D94").Value = -980500
$ws.Range("E94").Value = -640400
$ws.Range("F94").Value = -424200
$ws.Range("G94").Value = -490600
$ws.Range("H94").Value = -232200
$ws.Range("I94").Value = "NA"
$ws.Range("J94").Value = "NA"
$ws.Range("K94").Value = "NA"
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("D100").Value = 256100
$ws.Range("E100").Value = 52900
$ws.Range("F100").Value = 258600
$ws.Range("G100").Value = 360000
$ws.Range("H100").Value = -118100
$ws.Range("I100").Value = "NA"
$ws.Range("J100").Value = "NA"
$ws.Range("K100").Value = "NA"
$ws.Range("D101").Value = -18600
$ws.Range("E101").Value = 1700
$ws.Range("F101").Value = 3700
$ws.Range("G101").Value = -12200
$ws.Range("H101").Value = -6700
$ws.Range("I101").Value = "NA"
$ws.Range("J101").Value = "NA"
$ws.Range("K101").Value = "NA"
$ws.Range("D102").Value = 73800
$ws.Range("E102").Value = -12600
$ws.Range("F102").Value = 306300
$ws.Range("G102").Value = 167400
$ws.Range("H102").Value = -67900
$ws.Range("I102").Value = "NA"
$ws.Range("J102").Value = "NA"
$ws.Range("K102").Value = "NA"
